# rbnz_wholesale_rates.xlsx update
# - append 17 new monthly rows (436-452) to the "Data" sheet
# - bump the "as at" date on the "Table Description" sheet (B4)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# ---------------------------------------------------------------------------
# 1. Copy the per-column number formatting down into the new rows first, so
#    the new cells pick up the same styles as the rest of the table
#    (A -> date style, B:J/L:O -> numeric style, O -> 2dp style). Column K is
#    intentionally left untouched - it has no formatting/values in this part
#    of the table either. Column G is only populated (and therefore only
#    formatted) on rows 436 and 447-452 - rows 437-446 have no G cell at all.
# ---------------------------------------------------------------------------
$formatCols = @("A","B","C","D","E","F","H","I","J","L","M","N","O")
foreach ($col in $formatCols) {
    $src = $col + "400"
    $dst = $col + "436:" + $col + "452"
    $ws.Range($src).Copy()
    $ws.Range($dst).PasteSpecial(-4122)
}

$ws.Range("G400").Copy()
$ws.Range("G436").PasteSpecial(-4122)
$ws.Range("G447:G452").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. Write the new row values.
# ---------------------------------------------------------------------------
$rows = @(
    @{ Row=436; A=44165; B=0.25; C=0.22; D=0.27; E=0.27; F=0.27; G=0.21; H=0.2;  I=0.21; J=0.75; L=-0.53; M=-0.41; N=-0.24; O=0.02 },
    @{ Row=437; A=44196; B=0.25; C=0.23; D=0.26; E=0.26; F=0.26;            H=0.24; I=0.34; J=0.92; L=-0.51; M=-0.35; N=-0.09; O=0.16 },
    @{ Row=438; A=44227; B=0.25; C=0.24; D=0.26; E=0.27; F=0.28000000000000003; H=0.22; I=0.39; J=1.04; L=-0.62; M=-0.42; N=-0.14000000000000001; O=0.14000000000000001 },
    @{ Row=439; A=44255; B=0.25; C=0.24; D=0.26; E=0.27; F=0.28000000000000003; H=0.27; I=0.75; J=1.46; L=-0.73; M=-0.31; N=0.13; O=0.56999999999999995 },
    @{ Row=440; A=44286; B=0.25; C=0.23; D=0.26; E=0.28999999999999998; F=0.33; H=0.32; I=1.03; J=1.76; L=-0.86; M=-0.22; N=0.33; O=0.86 },
    @{ Row=441; A=44316; B=0.25; C=0.24; D=0.26; E=0.3;  F=0.34; H=0.25; I=0.88; J=1.68; L=-1;    M=-0.2;  N=0.44; O=0.94 },
    @{ Row=442; A=44347; B=0.25; C=0.22; D=0.27; E=0.31; F=0.35; H=0.3;  I=1.03; J=1.81; L=-0.99; M=-0.14000000000000001; N=0.43; O=0.87 },
    @{ Row=443; A=44377; B=0.25; C=0.23; D=0.27; E=0.3;  F=0.33; H=0.35; I=1.03; J=1.76; L=-0.89; M=-0.04; N=0.53; O=0.97 },
    @{ Row=444; A=44408; B=0.25; C=0.23; D=0.3;  E=0.35; F=0.4;  H=0.65; I=1.1299999999999999; J=1.59; L=-0.8; M=-0.1; N=0.33; O=0.79 },
    @{ Row=445; A=44439; B=0.25; C=0.21; D=0.39; E=0.47; F=0.54; H=0.9;  I=1.3;  J=1.65; L=-0.72; M=-0.19; N=0.2;  O=0.71 },
    @{ Row=446; A=44469; B=0.25; C=0.2;  D=0.38; E=0.48; F=0.56999999999999995; H=1.02; I=1.52; J=1.87; L=-0.57999999999999996; M=-0.07; N=0.31; O=0.78 },
    @{ Row=447; A=44500; B=0.5;  C=0.38; D=0.54; E=0.62; F=0.7;  G=1.33; H=1.1100000000000001; I=1.79; J=2.21; L=-0.53; M=0.16; N=0.47; O=0.78 },
    @{ Row=448; A=44530; B=0.75; C=0.54; D=0.75; E=0.79; F=0.83; G=1.58; H=2.02; I=2.33; J=2.57; L=-0.35; M=0.41; N=0.68; O=0.86 },
    @{ Row=449; A=44561; B=0.75; C=0.71; D=0.79; E=0.85; F=0.91; G=1.52; H=1.97; I=2.21; J=2.38; L=-0.5;  M=0.23; N=0.54; O=0.69 },
    @{ Row=450; A=44592; B=0.75; C=0.67; D=0.82; E=0.92; F=1.03; G=1.6;  H=2.0499999999999998; I=2.36; J=2.56; L=-0.33; M=0.42; N=0.72; O=0.85 },
    @{ Row=451; A=44620; B=1;    C=0.73; D=1;    E=1.1100000000000001; F=1.21; G=1.82; H=2.25; I=2.58; J=2.74; L=-0.14000000000000001; M=0.61; N=0.87; O=1 },
    @{ Row=452; A=44651; B=1;    C=0.94; D=1.1299999999999999; E=1.31; F=1.49; G=2.14; H=2.64; I=2.93; J=3.07; L=-0.17; M=0.56999999999999995; N=0.79; O=0.91 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    if ($r.ContainsKey("G")) { $ws.Cells.Item($row, 7).Value = $r.G }
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
}

# ---------------------------------------------------------------------------
# 3. Update the "as at" date on the Table Description sheet.
# ---------------------------------------------------------------------------
$tds = $wb.Worksheets.Item("Table Description")
$tds.Range("B4").Value = 44652
